# Rename "wire_transfers" sheet concept to "currency conversion to EUR":
# - update the "Wire transfer" comment text on the Fees sheet to the new,
#   clearer wording
# - widen the Comment column so the longer text still "best fits"
# - refresh the two dependent numbers that shifted because of the
#   corrected currency-conversion date (Foreign Currencies sheet + the
#   ELSTER summary roll-up)

$wb = $excel.ActiveWorkbook

# --- Fees sheet: reword the "Wire transfer" comment -------------------
$fees = $wb.Worksheets.Item("Fees")
$fees.Range("B6").Value = "Currency conversion or wire transfer"
$fees.Range("B11").Value = "Currency conversion or wire transfer"

# Widen column B (Comment) to fit the new, longer text.
$fees.Columns.Item(2).ColumnWidth = 32.333333

# --- Foreign Currencies sheet: updated conversion amount/gain ---------
$foreign = $wb.Worksheets.Item("Foreign Currencies")
$foreign.Range("B7").Value = 155
$foreign.Range("G7").Value = -10.64

# --- ELSTER - Summary sheet: updated currency gain/loss roll-up -------
$elster = $wb.Worksheets.Item("ELSTER - Summary")
$elster.Range("C7").Value = 67.67
